# Apply updated crypto price/volume figures (and a few re-ranked rows)
# from the commit. Each cell is written with a leading apostrophe so Excel
# keeps numeric-looking values ("0.998", "530.52", ...) stored as literal text,
# matching the original inline-string cells rather than coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''64.347.51'
$ws.Range("E2").Value = '''  -0.73%  '
$ws.Range("D3").Value = '''3.326.11'
$ws.Range("E3").Value = '''  +0.41%  '
$ws.Range("D4").Value = '''0.998'
$ws.Range("E4").Value = '''  -0.33%  '
$ws.Range("D5").Value = '''530.52'
$ws.Range("E5").Value = '''  -0.08%  '
$ws.Range("D6").Value = '''176.77'
$ws.Range("E6").Value = '''  -2.77%  '
$ws.Range("D7").Value = '''0.593'
$ws.Range("E7").Value = '''  -2.14%  '
$ws.Range("D8").Value = '''3.317.34'
$ws.Range("E8").Value = '''  +0.31%  '
$ws.Range("E9").Value = '''  -0.07%  '
$ws.Range("D10").Value = '''0.611'
$ws.Range("E10").Value = '''  -0.92%  '
$ws.Range("D11").Value = '''54.03'
$ws.Range("E11").Value = '''  -9.46%  '
$ws.Range("D12").Value = '''0.139'
$ws.Range("E12").Value = '''  +3.74%  '
$ws.Range("D13").Value = '''0.0000261'
$ws.Range("E13").Value = '''  -0.21%  '
$ws.Range("D14").Value = '''9.02'
$ws.Range("E14").Value = '''  -1.59%  '
$ws.Range("D15").Value = '''3.853.75'
$ws.Range("E15").Value = '''  +0.41%  '
$ws.Range("E16").Value = '''  +0.33%  '
$ws.Range("D17").Value = '''3.316.96'
$ws.Range("E17").Value = '''  +0.11%  '
$ws.Range("D18").Value = '''64.327.55'
$ws.Range("E18").Value = '''  -0.53%  '
$ws.Range("D19").Value = '''17.56'
$ws.Range("E19").Value = '''  -0.93%  '
$ws.Range("D20").Value = '''11.28'
$ws.Range("E20").Value = '''  +0.58%  '
$ws.Range("D21").Value = '''0.962'
$ws.Range("E21").Value = '''  -0.30%  '
$ws.Range("D22").Value = '''384.34'
$ws.Range("E22").Value = '''  +1.98%  '
$ws.Range("D23").Value = '''4.17'
$ws.Range("E23").Value = '''  +6.01%  '
$ws.Range("D24").Value = '''82.20'
$ws.Range("E24").Value = '''  +1.07%  '
$ws.Range("D25").Value = '''11.19'
$ws.Range("E25").Value = '''  +0.25%  '
$ws.Range("D26").Value = '''3.73'
$ws.Range("E26").Value = '''  -2.84%  '
$ws.Range("E27").Value = '''  -1.11%  '
$ws.Range("D28").Value = '''2.75'
$ws.Range("E28").Value = '''  +1.68%  '
$ws.Range("D29").Value = '''11.35'
$ws.Range("E29").Value = '''  -2.28%  '
$ws.Range("D30").Value = '''8.28'
$ws.Range("E30").Value = '''  -2.12%  '
$ws.Range("D31").Value = '''29.07'
$ws.Range("E31").Value = '''  -0.27%  '
$ws.Range("D32").Value = '''641.37'
$ws.Range("E32").Value = '''  -1.62%  '
$ws.Range("D33").Value = '''6.76'
$ws.Range("E33").Value = '''  +0.45%  '
$ws.Range("D34").Value = '''11.26'
$ws.Range("E34").Value = '''  -0.95%  '
$ws.Range("E35").Value = '''  -0.16%  '
$ws.Range("D36").Value = '''57.46'
$ws.Range("E36").Value = '''  -2.92%  '
$ws.Range("E37").Value = '''  -0.05%  '
$ws.Range("D38").Value = '''36.63'
$ws.Range("E38").Value = '''  -1.32%  '
$ws.Range("D39").Value = '''0.384'
$ws.Range("E39").Value = '''  -2.96%  '
$ws.Range("D40").Value = '''0.0₃0760'
$ws.Range("E40").Value = '''  +7.83%  '
$ws.Range("E41").Value = '''  -0.26%  '
$ws.Range("D42").Value = '''3.25'
$ws.Range("E42").Value = '''  +13.88%  '
$ws.Range("D43").Value = '''0.127'
$ws.Range("E43").Value = '''  -0.18%  '
$ws.Range("B44").Value = '''Maker'
$ws.Range("C44").Value = '''https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '''2.977.85'
$ws.Range("E44").Value = '''  +2.89%  '
$ws.Range("B45").Value = '''Fetch.AI'
$ws.Range("C45").Value = '''https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").Value = '''2.63'
$ws.Range("E45").Value = '''  +5.28%  '
$ws.Range("D46").Value = '''0.0403'
$ws.Range("E46").Value = '''  +0.22%  '
$ws.Range("D47").Value = '''2.69'
$ws.Range("E47").Value = '''  +0.88%  '
$ws.Range("B48").Value = '''ThetaToken'
$ws.Range("C48").Value = '''https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D48").Value = '''2.68'
$ws.Range("E48").Value = '''  -1.38%  '
$ws.Range("B49").Value = '''ApeXProtocol'
$ws.Range("C49").Value = '''https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").Value = '''3.11'
$ws.Range("E49").Value = '''  +1.08%  '
$ws.Range("B50").Value = '''Stellar'
$ws.Range("C50").Value = '''https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").Value = '''0.126'
$ws.Range("E50").Value = '''  -1.24%  '
$ws.Range("B51").Value = '''Monero'
$ws.Range("C51").Value = '''https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").Value = '''138.24'
$ws.Range("E51").Value = '''  +1.86%  '
